# The workbook gets a new data row inserted at row 415 (pushing the old
# rows 415..518 down to 416..519, which is why nearly every row below 415
# shows a "shift" in the diff). The new row carries the same
# Mercado/Region/Categoria/Variedad/Calidad/Unidad/Origen/Kg-Unidades/
# Clasificacion values as the (old) row 415, but with fresh date/volume/
# price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 415; Excel shifts rows 415:518 down to 416:519
# and copies formatting (incl. the date style on column D) from the row
# above, matching the style="2" seen on D415 in the diff.
$ws.Rows.Item(415).Insert()

$ws.Range("A415").Value = 10
$ws.Range("B415").Value = "Vega Modelo de Temuco"
$ws.Range("C415").Value = "La Araucanía"
$ws.Range("D415").Value = 44722
$ws.Range("E415").Value = 9
$ws.Range("F415").Value = 100112043
$ws.Range("G415").Value = "Pepino ensalada"
$ws.Range("H415").Value = "Sin especificar"
$ws.Range("I415").Value = "Primera"
$ws.Range("J415").Value = 80
$ws.Range("K415").Value = 22000
$ws.Range("L415").Value = 22000
$ws.Range("M415").Value = 22000
$ws.Range("N415").Value = "`$/caja 60 unidades"
$ws.Range("O415").Value = "Región de Arica y Parinacota"
$ws.Range("P415").Value = 367
$ws.Range("Q415").Value = 60
$ws.Range("R415").Value = "Hortaliza"
